$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update working_capital (column C) values
$ws.Range("C2").Value = 1000
$ws.Range("C3").Value = 1000
$ws.Range("C4").Value = 3000
$ws.Range("C5").Value = 3000
$ws.Range("C6").Value = 5000
$ws.Range("C7").Value = 5000

# Update fixed_cost (column K) values
$ws.Range("K2").Value = 50
$ws.Range("K3").Value = 50
$ws.Range("K4").Value = 26
$ws.Range("K5").Value = 26
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 5

# Update the active selection to K2:K3 with active cell K2
$ws.Range("K2:K3").Select()
